$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert "미래에셋비전스팩6호" as the new row 3 (pushes existing rows down) ---
$ws.Rows(3).Insert()
$ws.Range("A2:F2").Copy($ws.Range("A3:F3"))
$ws.Range("A3").Value = "미래에셋비전스팩6호"
$ws.Range("B3").Value = "2024.06.04~06.05"
$ws.Range("C3").Value = "2,000~2,000"
$ws.Range("D3").Value = "-"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "12900"
$ws.Range("F3").Value = "미래에셋증권"

# --- Insert "엑셀세라퓨틱스" as the new row 8 (pushes existing rows down) ---
$ws.Rows(8).Insert()
$ws.Range("A7:F7").Copy($ws.Range("A8:F8"))
$ws.Range("A8").Value = "엑셀세라퓨틱스"
$ws.Range("B8").Value = "2024.06.03~06.10"
$ws.Range("C8").Value = "6,200~7,700"
$ws.Range("D8").Value = "-"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "10032"
$ws.Range("F8").Value = "대신증권"

# Restore default (unstyled) appearance on the two forced-text cells so they
# match the plain (style-less) data cells around them.
$ws.Range("E3").Style = "Normal"
$ws.Range("E8").Style = "Normal"

# --- Remove the last two rows (SK증권스팩12호, HD현대마린솔루션...) ---
$ws.Rows(22).Delete()
$ws.Rows(22).Delete()
